$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Loss-of-sale walk-in records received for Kannur SG - 23 to 25 Dec 2025 intake.
# New rows are appended directly below the existing data (through row 40).

# Row 41
$ws.Cells.Item(41, 1).Value = 39
$ws.Cells.Item(41, 1).NumberFormat = "0"
$ws.Cells.Item(41, 2).Value = "23-12-2025"
$ws.Cells.Item(41, 3).Value = "ILAN"
$ws.Cells.Item(41, 4).Value = 8714310451
$ws.Cells.Item(41, 4).NumberFormat = "0"
$ws.Cells.Item(41, 5).Value = "24-12-2025"
$ws.Cells.Item(41, 6).Value = "ASWANTH. K"
$ws.Cells.Item(41, 7).Value = "Loss"
$ws.Cells.Item(41, 8).Value = "SIZE NOT SUITABLE"
$ws.Cells.Item(41, 9).Value = "SIZE TOO SMALL"
$ws.Cells.Item(41, 10).Value = "-"
$ws.Cells.Item(41, 11).Value = "LOSS"

# Row 42
$ws.Cells.Item(42, 1).Value = 40
$ws.Cells.Item(42, 1).NumberFormat = "0"
$ws.Cells.Item(42, 2).Value = "23-12-2025"
$ws.Cells.Item(42, 3).Value = "GEORGE"
$ws.Cells.Item(42, 4).Value = 8943368212
$ws.Cells.Item(42, 4).NumberFormat = "0"
$ws.Cells.Item(42, 5).Formula = "=""10-01-2026"""
$ws.Cells.Item(42, 5).Copy()
$ws.Cells.Item(42, 5).PasteSpecial(-4163)
$ws.Cells.Item(42, 6).Value = "AKHIL RAJ K"
$ws.Cells.Item(42, 7).Value = "Loss"
$ws.Cells.Item(42, 8).Value = "ENQUIRY"
$ws.Cells.Item(42, 9).Value = "ENQUIRY WITHOUT BRIDE/FAMILY"
$ws.Cells.Item(42, 10).Value = "-"
$ws.Cells.Item(42, 11).Value = "LOSS"

# Row 43
$ws.Cells.Item(43, 1).Value = 41
$ws.Cells.Item(43, 1).NumberFormat = "0"
$ws.Cells.Item(43, 2).Value = "24-12-2025"
$ws.Cells.Item(43, 3).Value = "navyuh"
$ws.Cells.Item(43, 4).Value = 8281535336
$ws.Cells.Item(43, 4).NumberFormat = "0"
$ws.Cells.Item(43, 5).Formula = "=""05-01-2026"""
$ws.Cells.Item(43, 5).Copy()
$ws.Cells.Item(43, 5).PasteSpecial(-4163)
$ws.Cells.Item(43, 6).Value = "SHIDHIN A V"
$ws.Cells.Item(43, 7).Value = "Loss"
$ws.Cells.Item(43, 8).Value = "ENQUIRY"
$ws.Cells.Item(43, 9).Value = "Enquiry for Relative/Friend"
$ws.Cells.Item(43, 10).Value = "-"
$ws.Cells.Item(43, 11).Value = "LOSS"

# Row 44
$ws.Cells.Item(44, 1).Value = 42
$ws.Cells.Item(44, 1).NumberFormat = "0"
$ws.Cells.Item(44, 2).Value = "24-12-2025"
$ws.Cells.Item(44, 3).Value = "amruth"
$ws.Cells.Item(44, 4).Value = 9207177946
$ws.Cells.Item(44, 4).NumberFormat = "0"
$ws.Cells.Item(44, 5).Formula = "=""06-04-2026"""
$ws.Cells.Item(44, 5).Copy()
$ws.Cells.Item(44, 5).PasteSpecial(-4163)
$ws.Cells.Item(44, 6).Value = "ASWANTH. K"
$ws.Cells.Item(44, 7).Value = "Loss"
$ws.Cells.Item(44, 8).Value = "PRODUCT"
$ws.Cells.Item(44, 9).Value = "REQUIRED DESIGN NOT AVAILABLE"
$ws.Cells.Item(44, 10).Value = "-"
$ws.Cells.Item(44, 11).Value = "LOSS"

# Row 45
$ws.Cells.Item(45, 1).Value = 43
$ws.Cells.Item(45, 1).NumberFormat = "0"
$ws.Cells.Item(45, 2).Value = "24-12-2025"
$ws.Cells.Item(45, 3).Value = "ARJUN"
$ws.Cells.Item(45, 4).Value = 8309706924
$ws.Cells.Item(45, 4).NumberFormat = "0"
$ws.Cells.Item(45, 5).Value = "14-02-2026"
$ws.Cells.Item(45, 6).Value = "SHIDHIN A V"
$ws.Cells.Item(45, 7).Value = "Loss"
$ws.Cells.Item(45, 8).Value = "PRICING"
$ws.Cells.Item(45, 9).Value = "RENT TO HIGH"
$ws.Cells.Item(45, 10).Value = "-"
$ws.Cells.Item(45, 11).Value = "LOSS"

# Row 46
$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 1).NumberFormat = "0"
$ws.Cells.Item(46, 2).Value = "24-12-2025"
$ws.Cells.Item(46, 3).Value = "STAMIL"
$ws.Cells.Item(46, 4).Value = 9446100871
$ws.Cells.Item(46, 4).NumberFormat = "0"
$ws.Cells.Item(46, 5).Value = "30-12-2025"
$ws.Cells.Item(46, 6).Value = "Thejus R"
$ws.Cells.Item(46, 7).Value = "Loss"
$ws.Cells.Item(46, 8).Value = "PRICING"
$ws.Cells.Item(46, 9).Value = "RENT TO HIGH"
$ws.Cells.Item(46, 10).Value = "-"
$ws.Cells.Item(46, 11).Value = "LOSS"

# Row 47
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 1).NumberFormat = "0"
$ws.Cells.Item(47, 2).Value = "25-12-2025"
$ws.Cells.Item(47, 3).Value = "NASEEF"
$ws.Cells.Item(47, 4).Value = 9526079736
$ws.Cells.Item(47, 4).NumberFormat = "0"
$ws.Cells.Item(47, 5).Value = "27-12-2025"
$ws.Cells.Item(47, 6).Value = "SHIDHIN A V"
$ws.Cells.Item(47, 7).Value = "Loss"
$ws.Cells.Item(47, 8).Value = "PRICING"
$ws.Cells.Item(47, 9).Value = "RENT TO HIGH"
$ws.Cells.Item(47, 10).Value = "-"
$ws.Cells.Item(47, 11).Value = "LOSS"

# Row 48
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 1).NumberFormat = "0"
$ws.Cells.Item(48, 2).Value = "25-12-2025"
$ws.Cells.Item(48, 3).Value = "FATAL"
$ws.Cells.Item(48, 4).Value = 9526271844
$ws.Cells.Item(48, 4).NumberFormat = "0"
$ws.Cells.Item(48, 5).Formula = "=""03-01-2026"""
$ws.Cells.Item(48, 5).Copy()
$ws.Cells.Item(48, 5).PasteSpecial(-4163)
$ws.Cells.Item(48, 6).Value = "SHIDHIN A V"
$ws.Cells.Item(48, 7).Value = "Loss"
$ws.Cells.Item(48, 8).Value = "ENQUIRY"
$ws.Cells.Item(48, 9).Value = "ENQUIRY WITHOUT BRIDE/FAMILY"
$ws.Cells.Item(48, 10).Value = "-"
$ws.Cells.Item(48, 11).Value = "LOSS"

# Row 49
$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 1).NumberFormat = "0"
$ws.Cells.Item(49, 2).Value = "25-12-2025"
$ws.Cells.Item(49, 3).Value = "AMEER"
$ws.Cells.Item(49, 4).Value = 7306823944
$ws.Cells.Item(49, 4).NumberFormat = "0"
$ws.Cells.Item(49, 5).Value = "27-12-2025"
$ws.Cells.Item(49, 6).Value = "Thejus R"
$ws.Cells.Item(49, 7).Value = "Loss"
$ws.Cells.Item(49, 8).Value = "ENQUIRY"
$ws.Cells.Item(49, 9).Value = "Enquiry for Relative/Friend"
$ws.Cells.Item(49, 10).Value = "-"
$ws.Cells.Item(49, 11).Value = "WILL COME"

# Row 50
$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(50, 1).NumberFormat = "0"
$ws.Cells.Item(50, 2).Value = "25-12-2025"
$ws.Cells.Item(50, 3).Value = "JITHIN"
$ws.Cells.Item(50, 4).Value = 9496028978
$ws.Cells.Item(50, 4).NumberFormat = "0"
$ws.Cells.Item(50, 5).Value = "22-02-2026"
$ws.Cells.Item(50, 6).Value = "FARIZ V A"
$ws.Cells.Item(50, 7).Value = "Loss"
$ws.Cells.Item(50, 8).Value = "ENQUIRY"
$ws.Cells.Item(50, 9).Value = "ENQUIRY WITHOUT BRIDE/FAMILY"
$ws.Cells.Item(50, 10).Value = "-"
$ws.Cells.Item(50, 11).Value = "WILL COME"

$excel.CutCopyMode = 0
Write-Output "Added rows 41-50 (records #39-#48) to the loss-of-sale log."